$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 198.17
$ws.Range("I15").Value = 198.17
$ws.Range("K15").Value = 594.51
$ws.Range("M15").Value = -425.51
$ws.Range("H109").Value = 27675
$ws.Range("J109").Value = 27675
$ws.Range("L109").Value = 27675
$ws.Range("N109").Value = -30449
$ws.Range("H116").Value = 7131.8945
$ws.Range("I116").Value = 1836.3636
$ws.Range("J116").Value = 14413.25
$ws.Range("K116").Value = 1836.3636
$ws.Range("L116").Value = 14413.25
$ws.Range("M116").Value = 1605.6364
$ws.Range("N116").Value = -21297.25
$ws.Range("H132").Value = 100561.02
$ws.Range("I132").Value = 126781.625
$ws.Range("J132").Value = 7332.222
$ws.Range("K132").Value = 380344.875
$ws.Range("L132").Value = 21996.666
$ws.Range("M132").Value = -377814.875
$ws.Range("N132").Value = -27056.666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 825.8333
$ws.Range("I2").Value = 859.0714
$ws.Range("K2").Value = 859.0714
$ws.Range("M2").Value = -746.0714
$ws.Range("H61").Value = 1664.6765
$ws.Range("I61").Value = 913.13043
$ws.Range("J61").Value = 3236.0908
$ws.Range("K61").Value = 913.13043
$ws.Range("L61").Value = 3236.0908
$ws.Range("M61").Value = -701.13043
$ws.Range("N61").Value = -3660.0908
$ws.Range("H116").Value = 825.8333
$ws.Range("I116").Value = 859.0714
$ws.Range("K116").Value = 859.0714
$ws.Range("M116").Value = 1434.9286
$ws.Range("H132").Value = 2982.372
$ws.Range("I132").Value = 2340.2173
$ws.Range("J132").Value = 3720.85
$ws.Range("K132").Value = 7020.651899999999
$ws.Range("L132").Value = 11162.55
$ws.Range("M132").Value = -4490.651899999999
$ws.Range("N132").Value = -16222.55
$ws.Range("H136").Value = 1664.6765
$ws.Range("I136").Value = 913.13043
$ws.Range("J136").Value = 3236.0908
$ws.Range("K136").Value = 2739.39129
$ws.Range("L136").Value = 9708.2724
$ws.Range("M136").Value = -189.39129
$ws.Range("N136").Value = -14808.2724
$ws.Range("H137").Value = 51750
$ws.Range("J137").Value = 51750
$ws.Range("L137").Value = 51750
$ws.Range("N137").Value = -61950

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 825.8333
$ws.Range("I3").Value = 859.0714
$ws.Range("K3").Value = 859.0714
$ws.Range("M3").Value = -745.0714
$ws.Range("H59").Value = 51800
$ws.Range("J59").Value = 51800
$ws.Range("L59").Value = 51800
$ws.Range("N59").Value = -53494
$ws.Range("H137").Value = 32967.5
$ws.Range("J137").Value = 32967.5
$ws.Range("L137").Value = 32967.5
$ws.Range("N137").Value = -43167.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3503.75
$ws.Range("I5").Value = 646
$ws.Range("J5").Value = 8266.666999999999
$ws.Range("K5").Value = 1938
$ws.Range("L5").Value = 24800.001
$ws.Range("M5").Value = -1826
$ws.Range("N5").Value = -25024.001
$ws.Range("H68").Value = 7964.143
$ws.Range("I68").Value = 700
$ws.Range("J68").Value = 26124.5
$ws.Range("K68").Value = 2100
$ws.Range("L68").Value = 78373.5
$ws.Range("M68").Value = -1289
$ws.Range("N68").Value = -79995.5
$ws.Range("H71").Value = 7964.143
$ws.Range("I71").Value = 700
$ws.Range("J71").Value = 26124.5
$ws.Range("K71").Value = 6300
$ws.Range("L71").Value = 235120.5
$ws.Range("M71").Value = -2244
$ws.Range("N71").Value = -243232.5
$ws.Range("H75").Value = 2906.5
$ws.Range("J75").Value = 5000
$ws.Range("L75").Value = 15000
$ws.Range("N75").Value = -16996
$ws.Range("H78").Value = 2906.5
$ws.Range("J78").Value = 5000
$ws.Range("L78").Value = 45000
$ws.Range("N78").Value = -54984
$ws.Range("H86").Value = 7284.3335
$ws.Range("I86").Value = 8271.429
$ws.Range("K86").Value = 24814.287
$ws.Range("M86").Value = -23628.287
$ws.Range("H89").Value = 7284.3335
$ws.Range("I89").Value = 8271.429
$ws.Range("K89").Value = 74442.861
$ws.Range("M89").Value = -68514.861
$ws.Range("H113").Value = 660.81134
$ws.Range("I113").Value = 581.9231
$ws.Range("J113").Value = 880.5714
$ws.Range("K113").Value = 1745.7693
$ws.Range("L113").Value = 2641.7142
$ws.Range("M113").Value = 424.2307000000001
$ws.Range("N113").Value = -6981.7142
$ws.Range("H132").Value = 2644.3
$ws.Range("I132").Value = 953.6667
$ws.Range("J132").Value = 3368.8572
$ws.Range("K132").Value = 8583.0003
$ws.Range("L132").Value = 30319.7148
$ws.Range("M132").Value = -6053.0003
$ws.Range("N132").Value = -35379.7148
$ws.Range("H135").Value = 3503.75
$ws.Range("I135").Value = 646
$ws.Range("J135").Value = 8266.666999999999
$ws.Range("K135").Value = 5814
$ws.Range("L135").Value = 74400.003
$ws.Range("M135").Value = -3279
$ws.Range("N135").Value = -79470.003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 32304
$ws.Range("J46").Value = 32304
$ws.Range("L46").Value = 32304
$ws.Range("N46").Value = -32616
$ws.Range("H132").Value = 3845.9048
$ws.Range("I132").Value = 2454.375
$ws.Range("J132").Value = 4702.231
$ws.Range("K132").Value = 7363.125
$ws.Range("L132").Value = 14106.693
$ws.Range("M132").Value = -4833.125
$ws.Range("N132").Value = -19166.693
$ws.Range("H134").Value = 50440.723
$ws.Range("I134").Value = 19296
$ws.Range("J134").Value = 52272.766
$ws.Range("K134").Value = 57888
$ws.Range("L134").Value = 156818.298
$ws.Range("M134").Value = -55353
$ws.Range("N134").Value = -161888.298
$ws.Range("H137").Value = 72719.78
$ws.Range("J137").Value = 72719.78
$ws.Range("L137").Value = 72719.78
$ws.Range("N137").Value = -82919.78

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2591
$ws.Range("I16").Value = 2306.3333
$ws.Range("J16").Value = 3160.3333
$ws.Range("K16").Value = 2306.3333
$ws.Range("L16").Value = 3160.3333
$ws.Range("M16").Value = -2136.3333
$ws.Range("N16").Value = -3500.3333
$ws.Range("H132").Value = 4365.769
$ws.Range("I132").Value = 1916.2941
$ws.Range("J132").Value = 5555.514
$ws.Range("K132").Value = 5748.8823
$ws.Range("L132").Value = 16666.542
$ws.Range("M132").Value = -3218.8823
$ws.Range("N132").Value = -21726.542
$ws.Range("H136").Value = 4694.037
$ws.Range("I136").Value = 2352.4375
$ws.Range("K136").Value = 7057.3125
$ws.Range("M136").Value = -4507.3125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 5849354.5
$ws.Range("I132").Value = 515.8857400000001
$ws.Range("K132").Value = 1547.65722
$ws.Range("M132").Value = 982.3427799999999
$ws.Range("H136").Value = 2696.7827
$ws.Range("I136").Value = 1319.0769
$ws.Range("J136").Value = 4487.8
$ws.Range("K136").Value = 3957.2307
$ws.Range("L136").Value = 13463.4
$ws.Range("M136").Value = -1407.2307
$ws.Range("N136").Value = -18563.4
